# [Kadastro App] Yeni kayit eklendi: 3017
# Appends the new record (row 76) to both the master "Kayitlar" sheet and
# the filtered "Erdemli" district sheet, mirroring the existing rows.

$wb = $excel.ActiveWorkbook

$sheetNames = @("Kayitlar", "Erdemli")

$newRow = 76
# Ordered column => value pairs for the new record.
$values = [ordered]@{
    "A" = "3017"
    "B" = "2025-09-11"
    "C" = "Erdemli"
    "D" = "1"
    "E" = "3B"
    "F" = "SERDAR ARSLAN (Tekniker), ÖZKAN AKBAŞ (Mühendis)"
}
# Columns that look numeric/date-like and must be kept as Text so they
# don't get auto-converted into Number/Date cells (matches the rest of
# the column, which stores e.g. "3016" and "2025-09-11" as plain text).
$textFormatCols = @("A", "B", "D")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    foreach ($col in $textFormatCols) {
        $ws.Range($col + $newRow).NumberFormat = "@"
    }

    foreach ($col in $values.Keys) {
        $ws.Range($col + $newRow).Value = $values[$col]
    }
}
